# A new weekly price record was inserted at row 141 of the "Albahaca" sheet,
# pushing the existing rows 141-237 down to 142-238.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("141:141").Insert()

$ws.Range("A141").Value = 6
$ws.Range("B141").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C141").Value = 'Metropolitana'
$ws.Range("D141").Value = 44438
$ws.Range("E141").Value = 13
$ws.Range("F141").Value = 100112052
$ws.Range("G141").Value = 'Albahaca'
$ws.Range("H141").Value = 'Sin especificar'
$ws.Range("I141").Value = 'Primera'
$ws.Range("J141").Value = 200
$ws.Range("K141").Value = 5000
$ws.Range("L141").Value = 5500
$ws.Range("M141").Value = 5300
$ws.Range("N141").Value = '$/paquete'
$ws.Range("O141").Value = 'Región de Arica y Parinacota'
$ws.Range("P141").Value = 5300
$ws.Range("Q141").Value = 1
$ws.Range("R141").Value = 'Hortaliza'
